$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "2016-03-17 22:35:27"
$wsZhCn.Range("H3").Value = "2016-03-17 22:35:46"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "2016-03-17 22:35:30"
$wsDeDe.Range("H3").Value = "2016-03-17 22:35:52"
